# Refresh cryptos list with the latest price/volume snapshot.
# (Row 15/16 coins - ShibaInu and WrappedBTC - also swap rank positions.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.691.38"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.510.07"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'575.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'166.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "2.508.72"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("D13").Value = "'4.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "2.969.81"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.564.79"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'24.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "2.519.69"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'11.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "'349.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'1.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'70.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").Value = "'3.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'8.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "2.646.70"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").Value = "0.0₃0894"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'7.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'459.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "'1.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'159.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "'18.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "'4.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -6.75%  "
$ws.Range("D47").Value = "'142.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "'0.520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'0.0735"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'0.580"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
